# Plan.xlsx update: "Update plan and design"
#
# - Rows 13-16 on Sheet1 get the same "Done/green" treatment already present
#   on rows 2-12: a new Status column (D) set to "Done", a new duplicate
#   "Actual MD" column (F) mirroring (or correcting) the Est MD column (E),
#   a newly-touched (empty) Notes column (I), and every touched cell in the
#   row recoloured with the green font used throughout the rest of the table.
# - The active-cell selection on Sheet1 moves from E27 to G23.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$green = 5287936          # RGB(0,176,80) == style fontId 2 (FF00B050)
$xlCenter = -4108          # xlHAlignCenter

# Values to backfill into column F (Actual MD) for rows 13-16; these mirror
# column E (Est MD) except row 14, which differs (".2MD" vs ".2 MD").
$colFValues = @{
    13 = ".3MD"
    14 = ".2MD"
    15 = ".1MD"
    16 = ".2 MD"
}

foreach ($row in 13..16) {
    # New "Status" cell (D): value "Done", centered + green (matches C/D style).
    $dCell = $ws.Cells.Item($row, 4)
    $dCell.Value = "Done"
    $dCell.HorizontalAlignment = $xlCenter
    $dCell.Font.Color = $green

    # New "Actual MD" cell (F).
    $fCell = $ws.Cells.Item($row, 6)
    $fCell.Value = $colFValues[$row]
    $fCell.Font.Color = $green

    # Recolor the existing populated cells in the row to the green font.
    foreach ($col in 1, 2, 3, 5, 7, 8) {
        $ws.Cells.Item($row, $col).Font.Color = $green
    }
    # C already has centered alignment from the column style; keep it explicit.
    $ws.Cells.Item($row, 3).HorizontalAlignment = $xlCenter

    # Newly-touched (still empty) "Notes" cell (I) just picks up the green font.
    $ws.Cells.Item($row, 9).Font.Color = $green
}

# Move the remembered selection on Sheet1 from E27 to G23.
$ws.Range("G23").Select()
